$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Year of Treatment" column (B) is removed entirely, shifting every
# subsequent column one place to the left (C->B, D->C, ... I->H).
$ws.Columns.Item(2).Delete()

# Rename the (now shifted) header row, appending the ".deja.deja.deja" suffix.
$ws.Range("B1").Value = "Inject.deja.deja.deja"
$ws.Range("C1").Value = "Smoke / inhale.deja.deja.deja"
$ws.Range("D1").Value = "Eat / drink.deja.deja.deja"
$ws.Range("E1").Value = "Sniff.deja.deja.deja"
$ws.Range("F1").Value = "Other.deja.deja.deja"
$ws.Range("G1").Value = "Not known / missing.deja.deja.deja"
$ws.Range("H1").Value = "Total.deja.deja.deja"
